$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/sheet title
$ws.Name = "Booking information of wutemey"

# Update row 2 (booking) data
# Force text formatting on columns that hold text-like values (dates & numeric-looking codes)
# so Excel keeps them as plain strings instead of auto-converting to date serials / numbers.
$ws.Range("F2:H2").NumberFormat = "@"

$ws.Range("A2").Value = "5778f857c4d839489044db6a62007906"
$ws.Range("B2").Value = "Olson-Roberts"
$ws.Range("C2").Value = "8 Hudson Center"
$ws.Range("D2").Value = "Wildrye"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "01/01/2020"
$ws.Range("G2").Value = "02/02/2020"
$ws.Range("H2").Value = "52901"

# Revert the temporary text formatting so the cells keep the workbook's default style
$ws.Range("F2:H2").ClearFormats()

# Remove row 3 entirely (the second booking entry)
$ws.Rows.Item(3).Delete()
